# Update the "fault angle" (theta) low/high range on the three "level*"
# sheets from (0, 90) to (-90, 180), and leave the selection / active cell
# on each touched sheet the way the author left it when they made the edit.

$wb = $excel.ActiveWorkbook

# --- level1 -----------------------------------------------------------
$ws1 = $wb.Worksheets.Item("level1")
$ws1.Activate()
$ws1.Range("I2").Value = -90
$ws1.Range("J2").Value = 180
$ws1.Range("I2:J2").Select()

# --- level2 -----------------------------------------------------------
$ws2 = $wb.Worksheets.Item("level2")
$ws2.Activate()
$ws2.Range("I2").Value = -90
$ws2.Range("J2").Value = 180
$ws2.Range("I2:J2").Select()

# --- level3 -----------------------------------------------------------
$ws3 = $wb.Worksheets.Item("level3")
$ws3.Activate()
$ws3.Range("I2").Value = -90
$ws3.Range("J2").Value = 180
$ws3.Range("D41").Select()

# --- fixed --------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("fixed")
$ws4.Activate()
$ws4.Range("D36").Select()

# Leave "level3" as the active sheet/tab, matching the workbook's
# activeTab / tabSelected state.
$ws3.Activate()
